$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 10.088846
$ws.Range("H2").Value = 30.266538
$ws.Range("I2").Value = 0.1151445838515654
$ws.Range("J2").Value = 0.1151445838515654
$ws.Range("M2").Value = 255.0443116666667
$ws.Range("N2").Value = 765.132935
$ws.Range("O2").Value = 0.863617428561108
$ws.Range("P2").Value = 0.8636174285611079
$ws.Range("Q2").Value = 2573.102783581003
$ws.Range("R2").Value = 23157.92505222903
$ws.Range("S2").Value = 0.09944086941862781
$ws.Range("T2").Value = 0.0994408694186278

$ws.Range("G3").Value = 10.088846
$ws.Range("H3").Value = 30.266538
$ws.Range("I3").Value = 0.1151445838515654
$ws.Range("J3").Value = 0.1151445838515654
$ws.Range("M3").Value = 0.8952453333333334
$ws.Range("N3").Value = 2.685736
$ws.Range("O3").Value = 0.003031431940796009
$ws.Range("P3").Value = 0.003031431940796009
$ws.Range("Q3").Value = 9.031992300218668
$ws.Range("R3").Value = 81.28793070196801
$ws.Range("S3").Value = 0.0003490529692972997
$ws.Range("T3").Value = 0.0003490529692972997

$ws.Range("G4").Value = 10.088846
$ws.Range("H4").Value = 30.266538
$ws.Range("I4").Value = 0.1151445838515654
$ws.Range("J4").Value = 0.1151445838515654
$ws.Range("M4").Value = 7.050555333333333
$ws.Range("N4").Value = 21.151666
$ws.Range("O4").Value = 0.02387421396349043
$ws.Range("P4").Value = 0.02387421396349043
$ws.Range("Q4").Value = 71.13196697247866
$ws.Range("R4").Value = 640.1877027523079
$ws.Range("S4").Value = 0.002748986431609338
$ws.Range("T4").Value = 0.002748986431609338

$ws.Range("G5").Value = 10.088846
$ws.Range("H5").Value = 30.266538
$ws.Range("I5").Value = 0.1151445838515654
$ws.Range("J5").Value = 0.1151445838515654
$ws.Range("M5").Value = 32.33082866666666
$ws.Range("N5").Value = 96.99248599999999
$ws.Range("O5").Value = 0.1094769255346056
$ws.Range("P5").Value = 0.1094769255346056
$ws.Range("Q5").Value = 326.1807514703853
$ws.Range("R5").Value = 2935.626763233467
$ws.Range("S5").Value = 0.01260567503203098
$ws.Range("T5").Value = 0.01260567503203098

$ws.Range("I6").Value = 0.4327250566572728
$ws.Range("J6").Value = 0.4327250566572729
$ws.Range("M6").Value = 255.0443116666667
$ws.Range("N6").Value = 765.132935
$ws.Range("O6").Value = 0.863617428561108
$ws.Range("P6").Value = 0.8636174285611079
$ws.Range("Q6").Value = 9669.981952824075
$ws.Range("R6").Value = 87029.83757541668
$ws.Range("S6").Value = 0.3737089007043137
$ws.Range("T6").Value = 0.3737089007043137

$ws.Range("I7").Value = 0.4327250566572728
$ws.Range("J7").Value = 0.4327250566572729
$ws.Range("M7").Value = 0.8952453333333334
$ws.Range("N7").Value = 2.685736
$ws.Range("O7").Value = 0.003031431940796009
$ws.Range("P7").Value = 0.003031431940796009
$ws.Range("Q7").Value = 33.94314564442312
$ws.Range("R7").Value = 305.4883107998081
$ws.Range("S7").Value = 0.001311776558333619
$ws.Range("T7").Value = 0.00131177655833362

$ws.Range("I8").Value = 0.4327250566572728
$ws.Range("J8").Value = 0.4327250566572729
$ws.Range("M8").Value = 7.050555333333333
$ws.Range("N8").Value = 21.151666
$ws.Range("O8").Value = 0.02387421396349043
$ws.Range("P8").Value = 0.02387421396349043
$ws.Range("Q8").Value = 267.3211662129831
$ws.Range("R8").Value = 2405.890495916848
$ws.Range("S8").Value = 0.01033097058999925
$ws.Range("T8").Value = 0.01033097058999925

$ws.Range("I9").Value = 0.4327250566572728
$ws.Range("J9").Value = 0.4327250566572729
$ws.Range("M9").Value = 32.33082866666666
$ws.Range("N9").Value = 96.99248599999999
$ws.Range("O9").Value = 0.1094769255346056
$ws.Range("P9").Value = 0.1094769255346056
$ws.Range("Q9").Value = 1225.820437568201
$ws.Range("R9").Value = 11032.38393811381
$ws.Range("S9").Value = 0.04737340880462627
$ws.Range("T9").Value = 0.04737340880462627

$ws.Range("G10").Value = 15.69885766666667
$ws.Range("H10").Value = 47.096573
$ws.Range("I10").Value = 0.1791719719949428
$ws.Range("J10").Value = 0.1791719719949428
$ws.Range("M10").Value = 255.0443116666667
$ws.Range("N10").Value = 765.132935
$ws.Range("O10").Value = 0.863617428561108
$ws.Range("P10").Value = 0.8636174285611079
$ws.Range("Q10").Value = 4003.904347547973
$ws.Range("R10").Value = 36035.13912793175
$ws.Range("S10").Value = 0.1547360377244954
$ws.Range("T10").Value = 0.1547360377244953

$ws.Range("G11").Value = 15.69885766666667
$ws.Range("H11").Value = 47.096573
$ws.Range("I11").Value = 0.1791719719949428
$ws.Range("J11").Value = 0.1791719719949428
$ws.Range("M11").Value = 0.8952453333333334
$ws.Range("N11").Value = 2.685736
$ws.Range("O11").Value = 0.003031431940796009
$ws.Range("P11").Value = 0.003031431940796009
$ws.Range("Q11").Value = 14.05432906474756
$ws.Range("R11").Value = 126.488961582728
$ws.Range("S11").Value = 0.0005431476388008775
$ws.Range("T11").Value = 0.0005431476388008775

$ws.Range("G12").Value = 15.69885766666667
$ws.Range("H12").Value = 47.096573
$ws.Range("I12").Value = 0.1791719719949428
$ws.Range("J12").Value = 0.1791719719949428
$ws.Range("M12").Value = 7.050555333333333
$ws.Range("N12").Value = 21.151666
$ws.Range("O12").Value = 0.02387421396349043
$ws.Range("P12").Value = 0.02387421396349043
$ws.Range("Q12").Value = 110.6856646489575
$ws.Range("R12").Value = 996.170981840618
$ws.Range("S12").Value = 0.00427758999566778
$ws.Range("T12").Value = 0.004277589995667779

$ws.Range("G13").Value = 15.69885766666667
$ws.Range("H13").Value = 47.096573
$ws.Range("I13").Value = 0.1791719719949428
$ws.Range("J13").Value = 0.1791719719949428
$ws.Range("M13").Value = 32.33082866666666
$ws.Range("N13").Value = 96.99248599999999
$ws.Range("O13").Value = 0.1094769255346056
$ws.Range("P13").Value = 0.1094769255346056
$ws.Range("Q13").Value = 507.5570774833864
$ws.Range("R13").Value = 4568.013697350477
$ws.Range("S13").Value = 0.0196151966359788
$ws.Range("T13").Value = 0.0196151966359788

$ws.Range("G14").Value = 23.91632366666667
$ws.Range("H14").Value = 71.748971
$ws.Range("I14").Value = 0.2729583874962189
$ws.Range("J14").Value = 0.2729583874962189
$ws.Range("M14").Value = 255.0443116666667
$ws.Range("N14").Value = 765.132935
$ws.Range("O14").Value = 0.863617428561108
$ws.Range("P14").Value = 0.8636174285611079
$ws.Range("Q14").Value = 6099.722307162209
$ws.Range("R14").Value = 54897.50076445988
$ws.Range("S14").Value = 0.2357316207136711
$ws.Range("T14").Value = 0.235731620713671

$ws.Range("G15").Value = 23.91632366666667
$ws.Range("H15").Value = 71.748971
$ws.Range("I15").Value = 0.2729583874962189
$ws.Range("J15").Value = 0.2729583874962189
$ws.Range("M15").Value = 0.8952453333333334
$ws.Range("N15").Value = 2.685736
$ws.Range("O15").Value = 0.003031431940796009
$ws.Range("P15").Value = 0.003031431940796009
$ws.Range("Q15").Value = 21.41097715307289
$ws.Range("R15").Value = 192.698794377656
$ws.Range("S15").Value = 0.0008274547743642119
$ws.Range("T15").Value = 0.0008274547743642119

$ws.Range("G16").Value = 23.91632366666667
$ws.Range("H16").Value = 71.748971
$ws.Range("I16").Value = 0.2729583874962189
$ws.Range("J16").Value = 0.2729583874962189
$ws.Range("M16").Value = 7.050555333333333
$ws.Range("N16").Value = 21.151666
$ws.Range("O16").Value = 0.02387421396349043
$ws.Range("P16").Value = 0.02387421396349043
$ws.Range("Q16").Value = 168.6233633817429
$ws.Range("R16").Value = 1517.610270435686
$ws.Range("S16").Value = 0.006516666946214062
$ws.Range("T16").Value = 0.006516666946214061

$ws.Range("G17").Value = 23.91632366666667
$ws.Range("H17").Value = 71.748971
$ws.Range("I17").Value = 0.2729583874962189
$ws.Range("J17").Value = 0.2729583874962189
$ws.Range("M17").Value = 32.33082866666666
$ws.Range("N17").Value = 96.99248599999999
$ws.Range("O17").Value = 0.1094769255346056
$ws.Range("P17").Value = 0.1094769255346056
$ws.Range("Q17").Value = 773.234562803545
$ws.Range("R17").Value = 6959.111065231905
$ws.Range("S17").Value = 0.02988264506196959
$ws.Range("T17").Value = 0.02988264506196959

